$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the D9:D14 total formulas as one block so Excel collapses them
# back into a single shared formula group (si="1"), same as it does whenever
# a formula is filled down over an existing range.
$ws.Range("D9:D14").Formula = "=SUM(C9*B9)"

# --- Add the new expenditure row (row 15): "Right angle mini usb cable" ---
$ws.Range("A15").Value = "Right angle mini usb cable"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 2.39

# D15 already carries the shared formula (si="0", ref D3:D19) inherited from
# the template rows below it - just make sure it is the expected formula.
$ws.Range("D15").Formula = "=SUM(C15*B15)"

# E15: purchase date (28 Apr 2019 == serial 43583), formatted like the other
# date cells (copy number format from E14 so it reuses the existing style).
$ws.Range("E15").Value = 43583
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)

# F15: purchase location
$ws.Range("F15").Value = "Amazon"

# G15: the "Click Here" hyperlink cell
$ws.Range("G15").Value = "Click Here"
$ws.Hyperlinks.Add($ws.Range("G15"), "https://www.amazon.co.uk/Right-Angle-Mini-USB-Cable/dp/B01N5N6J7K")
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# --- Update the selected cell shown when the workbook was last saved ---
$ws.Range("F16").Select()

$wb.Application.CutCopyMode = $false
